$d = $word.ActiveDocument

# --- Content edit -------------------------------------------------------
# "This game is awesome" -> the real game description.
$d.Content.Find.Execute(
    "This game is awesome", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The Game will be a VR Game where you are a wizard and have to fight other wizards and drones.",
    2)

# The lone paragraph (text run + trailing hidden "_GoBack" bookmark) is
# split into three paragraphs: the text itself, a new blank paragraph,
# and a paragraph that only carries the bookmark.
$para = $d.Paragraphs(1)
$splitPos = $para.Range.End - 1   # just before the paragraph mark (and the bookmark)

$break1 = $d.Range($splitPos, $splitPos)
$break1.Text = "`r"

$break2 = $d.Range($splitPos, $splitPos)
$break2.Text = "`r"
